$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for the specified rows
$ws.Range("F2").Value = -6
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -2
$ws.Range("F8").Value = -9
